$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '22.381.47'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = '1.560.57'
$ws.Range("E3").Value = '  -0.83%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.003'
$ws.Range("E5").Value = '  +0.17%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '285.92'
$ws.Range("E6").Value = '  -1.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3645'
$ws.Range("E7").Value = '  -2.89%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '48.71'
$ws.Range("E8").Value = '  -2.53%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3339'
$ws.Range("E9").Value = '  -2.14%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.123'
$ws.Range("E10").Value = '  -2.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07377'
$ws.Range("E11").Value = '  -2.86%  '

$ws.Range("E12").Value = '  +0.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.71'
$ws.Range("E13").Value = '  -3.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.903'
$ws.Range("E14").Value = '  -1.71%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.854'
$ws.Range("E15").Value = '  -1.70%  '

$ws.Range("D16").Value = '1.561.16'
$ws.Range("E16").Value = '  -0.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001095'
$ws.Range("E17").Value = '  -2.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.88'
$ws.Range("E18").Value = '  -2.56%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06745'
$ws.Range("E19").Value = '  +0.10%  '

$ws.Range("E20").Value = '  +0.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.270'
$ws.Range("E21").Value = '  -0.37%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.00'
$ws.Range("E22").Value = '  -2.84%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.92'
$ws.Range("E23").Value = '  -2.47%  '

$ws.Range("D24").Value = '22.380.80'
$ws.Range("E24").Value = '  -0.37%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.393'
$ws.Range("E25").Value = '  +3.01%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.519'
$ws.Range("E26").Value = '  -3.59%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '149.49'
$ws.Range("E27").Value = '  +0.61%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.55'
$ws.Range("E28").Value = '  -3.24%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.994'
$ws.Range("E29").Value = '  -0.03%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '123.06'
$ws.Range("E30").Value = '  -2.51%  '

$ws.Range("D31").Value = '1.737.27'
$ws.Range("E31").Value = '  -0.36%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.052'
$ws.Range("E32").Value = '  +0.73%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.079'
$ws.Range("E33").Value = '  -1.39%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.990'
$ws.Range("E34").Value = '  +0.54%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.507'
$ws.Range("E35").Value = '  -4.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08225'
$ws.Range("E36").Value = '  -2.89%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02376'
$ws.Range("E37").Value = '  -4.00%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.305'
$ws.Range("E38").Value = '  -5.98%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06368'
$ws.Range("E39").Value = '  -2.83%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.2208'
$ws.Range("E40").Value = '  -4.04%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.300'
$ws.Range("E41").Value = '  -4.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.09'
$ws.Range("E42").Value = '  -3.16%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("B43").Value = 'Frax'
$ws.Range("C43").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D43").Value = '1.002'
$ws.Range("E43").Value = '  +0.10%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.6030'
$ws.Range("E44").Value = '  -4.39%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '13.67'
$ws.Range("E45").Value = '  -3.13%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.762'
$ws.Range("E46").Value = '  -1.52%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5721'
$ws.Range("E47").Value = '  -2.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = '124.37'
$ws.Range("E48").Value = '  -5.02%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").Value = '2.000'
$ws.Range("E49").Value = '  -4.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.207'
$ws.Range("E50").Value = '  -1.67%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07217'
$ws.Range("E51").Value = '  -1.71%  '
